$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.740.13'
$ws.Range("E2").Value = '  -1.73%  '

# Row 3
$ws.Range("D3").Value = '3.385.11'
$ws.Range("E3").Value = '  -2.10%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.61'
$ws.Range("E5").Value = '  -2.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.85'
$ws.Range("E6").Value = '  -3.84%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").Value = '3.385.74'
$ws.Range("E8").Value = '  -2.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("E9").Value = '  -0.32%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.49'
$ws.Range("E10").Value = '  -2.49%  '

# Row 11
$ws.Range("E11").Value = '  -2.16%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.401'
$ws.Range("E12").Value = '  +2.19%  '

# Row 13
$ws.Range("D13").Value = '3.962.99'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.44'
$ws.Range("E14").Value = '  +1.58%  '

# Row 15
$ws.Range("E15").Value = '  +1.54%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  -2.17%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.381.23'
$ws.Range("E17").Value = '  -2.22%  '

# Row 18
$ws.Range("D18").Value = '60.818.88'
$ws.Range("E18").Value = '  -1.75%  '

# Row 19
$ws.Range("E19").Value = '  +0.25%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.08'
$ws.Range("E20").Value = '  -2.25%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.02'
$ws.Range("E21").Value = '  -5.78%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.56'
$ws.Range("E22").Value = '  -1.32%  '

# Row 23
$ws.Range("E23").Value = '  -0.83%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.69'
$ws.Range("E24").Value = '  +0.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.21%  '

# Row 26
$ws.Range("E26").Value = '  -5.63%  '

# Row 27
$ws.Range("D27").Value = '3.523.44'
$ws.Range("E27").Value = '  -2.17%  '

# Row 28
$ws.Range("E28").Value = '  -2.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.44'
$ws.Range("E30").Value = '  -3.54%  '

# Row 31
$ws.Range("E31").Value = '  -2.16%  '

# Row 32
$ws.Range("E32").Value = '  -2.08%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.43'
$ws.Range("E33").Value = '  -3.47%  '

# Row 34
$ws.Range("E34").Value = '  -0.04%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.68'
$ws.Range("E35").Value = '  -2.26%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.01'
$ws.Range("E36").Value = '  -0.13%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.67'
$ws.Range("E37").Value = '  -0.24%  '

# Row 38
$ws.Range("E38").Value = '  -2.53%  '

# Row 39
$ws.Range("D39").Value = '3.416.46'
$ws.Range("E39").Value = '  -2.02%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.49'
$ws.Range("E40").Value = '  -4.97%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0778'
$ws.Range("E41").Value = '  -0.82%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.77'
$ws.Range("E42").Value = '  +1.33%  '

# Row 43
$ws.Range("E43").Value = '  -2.79%  '

# Row 44
$ws.Range("E44").Value = '  -0.08%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.44'
$ws.Range("E45").Value = '  -1.62%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.76'
$ws.Range("E46").Value = '  -2.09%  '

# Row 47
$ws.Range("E47").Value = '  -3.02%  '

# Row 48
$ws.Range("D48").Value = '2.524.56'
$ws.Range("E48").Value = '  -1.99%  '

# Row 49
$ws.Range("E49").Value = '  -4.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.64'
$ws.Range("E50").Value = '  +2.33%  '

# Row 51
$ws.Range("E51").Value = '  -1.00%  '
